# Apply crypto price/volume updates from the commit diff.
# Values are written with a leading apostrophe (text-literal prefix) so Excel
# preserves them as text (matching the original t="inlineStr" cells) instead of
# auto-converting number-looking strings (e.g. "31.08") into numeric values.
# The Style reset afterwards clears the transient quote-prefix formatting so no
# extra cell style is left applied.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''34.015.98'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -1.61%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''1.790.40'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -1.77%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  +0.28%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Value = '''  -2.24%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''0.554'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  +0.42%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '''  +0.07%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''31.08'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  -3.25%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''46.04'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  -0.55%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = '''  -1.92%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = '''  -3.52%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''0.0926'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  -0.58%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''2.048.56'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  -1.86%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''11.39'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  +9.80%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''1.793.16'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  -1.62%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''0.634'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -2.01%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''34.032.07'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  -1.53%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = '''  -2.80%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''69.41'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  -2.51%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''252.77'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  -3.75%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''0.0₃0743'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -2.50%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = '''  +0.45%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''10.42'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  -1.35%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = '''  -3.15%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = '''  -2.58%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''157.13'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  -3.17%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = '''  -2.33%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = '''  -2.54%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = '''  -2.89%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = '''  +0.30%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''3.87'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  +0.34%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = '''0.0516'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  -0.43%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = '''  -0.91%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = '''  +0.74%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = '''  -0.68%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = '''1.473.08'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  -7.33%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = '''  -0.39%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''0.631'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  -0.43%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = '''Aave'
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = '''https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = '''83.78'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  -2.70%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = '''VeChain'
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = '''https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = '''0.0186'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  -1.58%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = '''  +0.00%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = '''  -0.29%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = '''  -2.42%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = '''  -4.58%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = '''  -2.21%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = '''  -0.09%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''1.947.87'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  -1.57%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = '''  +0.29%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = '''  -0.56%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''11.67'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  +2.32%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''51.28'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  -5.25%  '
$ws.Range("E51").Style = "Normal"
